$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Issued this {{ ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Issued this{{ ",
    2
)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(
    "current_date }} at NIA-Pangasinan IMO, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "current_date }} at NIA-Pangasinan IMO, ",
    2
)

$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Execute(
    "Bayaoas, Urdaneta City.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bayaoas, Urdaneta City.",
    2
)
